$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column B (id_building_type) from its original best-fit width to
# accommodate the newly added cooling/ventilation technology labels.
$ws.Columns(2).ColumnWidth = 18.33

# Move the active selection/scroll position: the sheet view now shows cell
# C6 as the active cell (scrolled back to the top-left of the sheet) instead
# of the previous view that was scrolled down with H13 selected.
[void]$ws.Range("C6").Select()
